$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws1.Cells.Item(3, 2).Value = "6.0.0"

# Update Date value (row 8, column B)
$ws1.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (row 9, column B)
$ws1.Cells.Item(9, 2).Value = "Alvearie Team"

# Replace duplicate "Contact" row (row 10) with "Jurisdiction" / "United States of America"
$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"

# Delete the second duplicate "Contact" row (row 11), shifting rows 12-15 up
$ws1.Rows.Item(11).Delete()
